# Weekly fruit/vegetable price update: insert a new weekly record at row 24
# (shifting the existing historical rows down by one) and populate it with
# the new week's values, per the commit "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 24..139 down to 25..140, creating a blank row 24 and extending
# the sheet's used range / dimension to A1:R140 automatically.
$ws.Rows.Item(24).Insert()

# Copy the (now shifted-down) row 25 into the new blank row 24 so every
# "constant" column (Mercado, Region, Categoria, Unidad, Origen, etc.) is
# populated exactly like the rest of the data set.
$ws.Rows.Item(25).Copy()
$ws.Rows.Item(24).PasteSpecial()

# Now overwrite just the columns that hold this week's new values.
$ws.Range("D24").Value2 = 44687
$ws.Range("J24").Value2 = 45
$ws.Range("K24").Value2 = 9000
$ws.Range("L24").Value2 = 9000
$ws.Range("M24").Value2 = 9000
$ws.Range("P24").Value2 = 3000
